# Remove job posting #5 (Service Now Developer / Chennai / long description)
# from the jobs sheet, leaving row 6 with just the wrap-text styled, empty
# D6 cell and no explicit row height (matches Excel's "clear + autofit").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:D6").ClearContents()
$ws.Rows.Item(6).AutoFit()
